# Generate Report for Handoff
#
# Semantics of this update (localization-status.xlsx):
#   - File 904e6585-bd2f-40ee-bcea-711e1c3a430d.md moved from status
#     "In Translation" to "Ready for handoff".
#   - Because rows are listed alphabetically within a status group, it now
#     sorts before 2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md (which was
#     already "Ready for handoff"), so the two rows swap places (row 9/10
#     on every sheet).
#   - The "Latest Handoff Datetime" / "Latest Handoff Date" timestamps for
#     the affected rows (row 6 and rows 9-10) are bumped forward to reflect
#     the new report generation run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn), C (de-de),
# D (Latest Handoff Date). Column A is hyperlinked -> use TextToDisplay
# to keep the existing link target untouched, exactly like the diff.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 6 (bb00dec1...): status text unchanged, just the handoff date bumps.
$wsOverview.Cells.Item(6, 4).Value2 = "2016-03-24 11:10:24"

# Row 9 becomes 2e4c8da1 / Ready for handoff / Ready for handoff / new date
$ovHyperlinks = $wsOverview.Hyperlinks
$ovHyperlinks.Item(8).TextToDisplay = "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md"
$wsOverview.Cells.Item(9, 2).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(9, 3).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(9, 4).Value2 = "2016-03-24 11:10:24"

# Row 10 becomes 904e6585 / Ready for handoff / Ready for handoff / new date
$ovHyperlinks.Item(9).TextToDisplay = "904e6585-bd2f-40ee-bcea-711e1c3a430d.md"
$wsOverview.Cells.Item(10, 2).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(10, 3).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(10, 4).Value2 = "2016-03-24 11:10:24"

# ---------------------------------------------------------------------
# zh-cn sheet: A (Source File Name) & D (Latest Handoff File) are
# hyperlinked; C (Status) & E (Latest Handoff Datetime) are plain text.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 6: only the handoff datetime is bumped.
$wsZh.Cells.Item(6, 5).Value2 = "2016-03-24 11:10:19"

$zhHyperlinks = $wsZh.Hyperlinks

# Row 9 becomes 2e4c8da1 / Ready for handoff / new date
$zhHyperlinks.Item(23).TextToDisplay = "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md"
$zhHyperlinks.Item(24).TextToDisplay = "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.zh-cn.xlf"
$wsZh.Cells.Item(9, 3).Value2 = "Ready for handoff"
$wsZh.Cells.Item(9, 5).Value2 = "2016-03-24 11:10:19"

# Row 10 becomes 904e6585 / Ready for handoff / new date
$zhHyperlinks.Item(25).TextToDisplay = "904e6585-bd2f-40ee-bcea-711e1c3a430d.md"
$zhHyperlinks.Item(26).TextToDisplay = "904e6585-bd2f-40ee-bcea-711e1c3a430d.0f8e0222e820725feba4bb99759bf7fa0fa65d92.zh-cn.xlf"
$wsZh.Cells.Item(10, 3).Value2 = "Ready for handoff"
$wsZh.Cells.Item(10, 5).Value2 = "2016-03-24 11:10:19"

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 6: only the handoff datetime is bumped.
$wsDe.Cells.Item(6, 5).Value2 = "2016-03-24 11:10:24"

$deHyperlinks = $wsDe.Hyperlinks

# Row 9 becomes 2e4c8da1 / Ready for handoff / new date
$deHyperlinks.Item(23).TextToDisplay = "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md"
$deHyperlinks.Item(24).TextToDisplay = "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.de-de.xlf"
$wsDe.Cells.Item(9, 3).Value2 = "Ready for handoff"
$wsDe.Cells.Item(9, 5).Value2 = "2016-03-24 11:10:24"

# Row 10 becomes 904e6585 / Ready for handoff / new date
$deHyperlinks.Item(25).TextToDisplay = "904e6585-bd2f-40ee-bcea-711e1c3a430d.md"
$deHyperlinks.Item(26).TextToDisplay = "904e6585-bd2f-40ee-bcea-711e1c3a430d.0f8e0222e820725feba4bb99759bf7fa0fa65d92.de-de.xlf"
$wsDe.Cells.Item(10, 3).Value2 = "Ready for handoff"
$wsDe.Cells.Item(10, 5).Value2 = "2016-03-24 11:10:24"
